# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.629.85"
$ws.Range("E2").Value = "  -0.55%  "

# Row 3
$ws.Range("D3").Value = "2.428.38"
$ws.Range("E3").Value = "  -1.73%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'558.84"
$ws.Range("E5").Value = "  -0.27%  "

# Row 6
$ws.Range("D6").Value = "'160.38"
$ws.Range("E6").Value = "  -1.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  +0.50%  "

# Row 9
$ws.Range("D9").Value = "'0.166"
$ws.Range("E9").Value = "  +9.66%  "

# Row 10
$ws.Range("E10").Value = "  -1.52%  "

# Row 11
$ws.Range("D11").Value = "'0.330"
$ws.Range("E11").Value = "  -0.48%  "

# Row 12
$ws.Range("D12").Value = "'4.60"
$ws.Range("E12").Value = "  -5.57%  "

# Row 13
$ws.Range("D13").Value = "68.522.05"
$ws.Range("E13").Value = "  -0.46%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.877.12"
$ws.Range("E14").Value = "  -0.99%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000174"
$ws.Range("E15").Value = "  +3.36%  "

# Row 16
$ws.Range("D16").Value = "'23.11"
$ws.Range("E16").Value = "  -1.97%  "

# Row 17
$ws.Range("D17").Value = "2.431.91"
$ws.Range("E17").Value = "  -1.36%  "

# Row 18
$ws.Range("D18").Value = "'10.48"
$ws.Range("E18").Value = "  -2.18%  "

# Row 19
$ws.Range("D19").Value = "'334.92"
$ws.Range("E19").Value = "  -0.37%  "

# Row 20
$ws.Range("D20").Value = "'6.91"
$ws.Range("E20").Value = "  -0.79%  "

# Row 21
$ws.Range("D21").Value = "'3.82"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22
$ws.Range("D22").Value = "'1.93"
$ws.Range("E22").Value = "  +2.26%  "

# Row 24
$ws.Range("D24").Value = "'66.85"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25
$ws.Range("D25").Value = "'3.66"
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").Value = "2.555.67"
$ws.Range("E26").Value = "  -1.74%  "

# Row 27
$ws.Range("D27").Value = "'1.01"
$ws.Range("E27").Value = "  +1.21%  "

# Row 28
$ws.Range("D28").Value = "'8.19"
$ws.Range("E28").Value = "  +0.16%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0818"
$ws.Range("E29").Value = "  +0.19%  "

# Row 30
$ws.Range("D30").Value = "'7.14"
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.09%  "

# Row 32
$ws.Range("D32").Value = "'426.25"
$ws.Range("E32").Value = "  -0.90%  "

# Row 33
$ws.Range("D33").Value = "'1.14"
$ws.Range("E33").Value = "  +0.54%  "

# Row 34
$ws.Range("E34").Value = "  -0.26%  "

# Row 35
$ws.Range("D35").Value = "'160.67"
$ws.Range("E35").Value = "  +1.45%  "

# Row 36
$ws.Range("D36").Value = "'19.00"
$ws.Range("E36").Value = "  -0.13%  "

# Row 37
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("D38").Value = "'17.91"
$ws.Range("E38").Value = "  +0.80%  "

# Row 39
$ws.Range("E39").Value = "  -3.00%  "

# Row 40
$ws.Range("E40").Value = "  -0.64%  "

# Row 41
$ws.Range("E41").Value = "  -1.98%  "

# Row 42
$ws.Range("E42").Value = "  +1.81%  "

# Row 43
$ws.Range("D43").Value = "'1.07"
$ws.Range("E43").Value = "  +0.37%  "

# Row 44
$ws.Range("D44").Value = "'2.04"
$ws.Range("E44").Value = "  -1.06%  "

# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'3.34"
$ws.Range("E45").Value = "  -0.19%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'131.19"
$ws.Range("E46").Value = "  -0.69%  "

# Row 47
$ws.Range("D47").Value = "'0.0716"
$ws.Range("E47").Value = "  +0.41%  "

# Row 48
$ws.Range("D48").Value = "'0.481"
$ws.Range("E48").Value = "  -0.43%  "

# Row 49
$ws.Range("E49").Value = "  -0.76%  "

# Row 50
$ws.Range("E50").Value = "  +0.51%  "

# Row 51
$ws.Range("E51").Value = "  +0.24%  "
